$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.334.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.24%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.735.32'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.54%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.45%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.76%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.52%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4249'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -11.64%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3581'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.77'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.114'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.79%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07330'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.70%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.42%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.42'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.22%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.055'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.96%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.158'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.740.04'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.21%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001055'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.87%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '83.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05947'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -11.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.53%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.68'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.11%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.984'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.68%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.346.56'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.13%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.23'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.58%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.28%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.78'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.89%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.33%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.315'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.99%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.940.23'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.15%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.270'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.10%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.33'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.725'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.74%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08986'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.23%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.519'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.11%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.22'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.40%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2138'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.47%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02255'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.96%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06072'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.97%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6394'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.28%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.961'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.20%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.178'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.50%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.005'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.415'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.80%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.808'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.78%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.45'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -5.37%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.737'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.41%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5845'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.10%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.99'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.80%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.928'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.85%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.090'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.97%  '

